$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Songs")
$ws1.Range("B2").ClearContents()
$ws1.Activate()
$ws1.Range("B2").Select() | Out-Null
